$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Tgfb2", "Tgfbr1", "ECs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 63.76294533333333, 191.288836, 0.6446527016991613, 0.6446527016991614, 113.966636578588, 1025.699729207292, 0.03175244748078916, 0.03175244748078916),
    @("FAPs", "Tgfb2", "Tgfbr1", "FAPs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 22.32219333333333, 66.96658, 0.2256806388876402, 0.2256806388876402, 39.89754993214, 359.07794938926, 0.01111592740528813, 0.01111592740528813),
    @("FAPs", "Tgfb2", "Tgfbr1", "sCs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 12.825399, 38.476197, 0.1296666594131984, 0.1296666594131984, 22.923464077251, 206.311176695259, 0.006386747130935531, 0.006386747130935531),
    @("sCs", "Tgfb2", "Tgfbr1", "ECs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 63.76294533333333, 191.288836, 0.6446527016991613, 0.6446527016991614, 1241.453705939294, 11173.08335345364, 0.3458836268321929, 0.345883626832193),
    @("sCs", "Tgfb2", "Tgfbr1", "FAPs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 22.32219333333333, 66.96658, 0.2256806388876402, 0.2256806388876402, 434.6093094271333, 3911.4837848442, 0.1210872733155645, 0.1210872733155645),
    @("sCs", "Tgfb2", "Tgfbr1", "sCs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 12.825399, 38.476197, 0.1296666594131984, 0.1296666594131984, 249.70833821217, 2247.37504390953, 0.06957168459674218, 0.06957168459674218),
    @("ECs", "Tgfb2", "Tgfbr1", "ECs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 63.76294533333333, 191.288836, 0.6446527016991613, 0.6446527016991614, 958.3824035036124, 8625.441631532512, 0.2670166273861793, 0.2670166273861793),
    @("ECs", "Tgfb2", "Tgfbr1", "FAPs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 22.32219333333333, 66.96658, 0.2256806388876402, 0.2256806388876402, 335.5114351514844, 3019.60291636336, 0.09347743816678754, 0.09347743816678754),
    @("ECs", "Tgfb2", "Tgfbr1", "sCs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 12.825399, 38.476197, 0.1296666594131984, 0.1296666594131984, 192.770842928536, 1734.937586356824, 0.05370822768552069, 0.05370822768552069)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

Write-Output "done"
